# Add the new DDL test case "ddl_032" as row 33 of Sheet1, following the
# same layout as the existing rows (TestID, Testable, Title, Component,
# Sub_component, Table_schema_ref, Table_value_ref, Ddl_sql, Query_sql,
# Query_result, Validation_type).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 33

$ws.Cells.Item($row, 1).Value = "ddl_032"
$ws.Cells.Item($row, 2).Value = "y"
$ws.Cells.Item($row, 3).Value = "创建Database，再次使用if not exists创建相同database"
$ws.Cells.Item($row, 4).Value = "Schema"
$ws.Cells.Item($row, 5).Value = "Information_Schema"
$ws.Cells.Item($row, 8).Value = "create database MYDDL_032;create database if not exists MYDDL_032"
$ws.Cells.Item($row, 9).Value = "select * from information_schema.schemata where schema_name in ('MYDDL_032')"
$ws.Cells.Item($row, 10).Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/ddl/expectedresult/ddl_032.csv"
$ws.Cells.Item($row, 11).Value = "csv_containsAll"

# Match the saved view state: sheet scrolled/selected around the new row.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
[void]$ws.Range("G40").Select()
